$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the '/' separators with '-' in the date strings for rows 3-21 (column A),
# forcing the cell to keep its original "General" style/text type instead of letting
# Excel auto-detect the new text as a date value.
for ($r = 3; $r -le 21; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $date = $cell.Value()
    $cell.NumberFormat = "@"
    $cell.Value = $date.Replace("/", "-")
    $cell.Style = "Normal"
}

# Update the attendance counts that changed for rows 3 and 4
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0
